$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp shown in A1
$ws.Range("A1").Value = "Datos actualizados a 6 de Abril de 2020 a las 11:22"

# Updated country statistics (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes) for the rows whose
# figures changed in this refresh. Column order: B,C,D,E,F,G,H

# Row 14 - Belgica
$ws.Range("B14").Value = 20814
$ws.Range("C14").Value = 1123
$ws.Range("D14").Value = 3986
$ws.Range("E14").Value = 15196
$ws.Range("F14").Value = 1257
$ws.Range("G14").Value = 185
$ws.Range("H14").Value = 1632

# Row 25 - Noruega
$ws.Range("B25").Value = 5760
$ws.Range("C25").Value = 73
$ws.Range("D25").Value = 32
$ws.Range("E25").Value = 5655
$ws.Range("F25").Value = 89
$ws.Range("G25").Value = 2
$ws.Range("H25").Value = 73

# Row 63 - Eslovenia
$ws.Range("B63").Value = 1021
$ws.Range("C63").Value = 24
$ws.Range("D63").Value = 79
$ws.Range("E63").Value = 912
$ws.Range("F63").Value = 30
$ws.Range("G63").Value = 2
$ws.Range("H63").Value = 30

# Row 74 - Camerun
$ws.Range("B74").Value = 665
$ws.Range("C74").Value = 109
$ws.Range("D74").Value = 103
$ws.Range("E74").Value = 561
$ws.Range("F74").Value = 17
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 1

# Row 75 - Azerbaiyan
$ws.Range("B75").Value = 658
$ws.Range("C75").Value = 8
$ws.Range("D75").Value = 17
$ws.Range("E75").Value = 632
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 9

# Row 76 - Kazajistan
$ws.Range("B76").Value = 604
$ws.Range("C76").Value = 20
$ws.Range("D76").Value = 45
$ws.Range("E76").Value = 553
$ws.Range("F76").Value = 6
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 6

# Row 77 - Tunez
$ws.Range("B77").Value = 584
$ws.Range("C77").Value = 0
$ws.Range("D77").Value = 32
$ws.Range("E77").Value = 545
$ws.Range("F77").Value = 17
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 7

# Row 78 - Kuwait
$ws.Range("B78").Value = 574
$ws.Range("C78").Value = 0
$ws.Range("D78").Value = 5
$ws.Range("E78").Value = 547
$ws.Range("F78").Value = 39
$ws.Range("G78").Value = 0
$ws.Range("H78").Value = 22

# Row 100 - Estado de Palestina
$ws.Range("B100").Value = 252
$ws.Range("C100").Value = 15
$ws.Range("D100").Value = 25
$ws.Range("E100").Value = 226
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 1

# Row 119 - Guadalupe
$ws.Range("B119").Value = 138
$ws.Range("C119").Value = 11
$ws.Range("D119").Value = 0
$ws.Range("E119").Value = 137
$ws.Range("F119").Value = 0
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 1

# Row 120 - Brunei
$ws.Range("B120").Value = 135
$ws.Range("C120").Value = 0
$ws.Range("D120").Value = 31
$ws.Range("E120").Value = 97
$ws.Range("F120").Value = 14
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 7

# Row 121 - Isla de Man
$ws.Range("B121").Value = 135
$ws.Range("C121").Value = 0
$ws.Range("D121").Value = 73
$ws.Range("E121").Value = 61
$ws.Range("F121").Value = 3
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 1

# Row 122 - Banglades
$ws.Range("B122").Value = 123
$ws.Range("C122").Value = 35
$ws.Range("D122").Value = 33
$ws.Range("E122").Value = 78
$ws.Range("F122").Value = 1
$ws.Range("G122").Value = 3
$ws.Range("H122").Value = 12

# Row 187 - Zimbabue
$ws.Range("B187").Value = 9
$ws.Range("C187").Value = 0
$ws.Range("D187").Value = 1
$ws.Range("E187").Value = 8
$ws.Range("F187").Value = 0
$ws.Range("G187").Value = 0
$ws.Range("H187").Value = 0

# Row 188 - Nepal
$ws.Range("B188").Value = 9
$ws.Range("C188").Value = 0
$ws.Range("D188").Value = 0
$ws.Range("E188").Value = 8
$ws.Range("F188").Value = 0
$ws.Range("G188").Value = 0
$ws.Range("H188").Value = 1

# Row 193 - Cabo Verde
$ws.Range("B193").Value = 7
$ws.Range("C193").Value = 0
$ws.Range("D193").Value = 1
$ws.Range("E193").Value = 6
$ws.Range("F193").Value = 0
$ws.Range("G193").Value = 0
$ws.Range("H193").Value = 0

# Row 194 - Somalia
$ws.Range("B194").Value = 7
$ws.Range("C194").Value = 0
$ws.Range("D194").Value = 1
$ws.Range("E194").Value = 5
$ws.Range("F194").Value = 0
$ws.Range("G194").Value = 0
$ws.Range("H194").Value = 1
